$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: add Position value (D column)
$ws.Range("D31").Value = 17000

# Row 40
$ws.Range("A40").Value = 0.09
$ws.Range("B40").NumberFormat = "0%"
$ws.Range("B40").Value = 0.96
$ws.Range("C40").NumberFormat = "0.00%"
$ws.Range("C40").Value = 0.0112

# Row 41
$ws.Range("A41").Value = 0.2
$ws.Range("B41").NumberFormat = "0%"
$ws.Range("B41").Value = 0.95
$ws.Range("C41").NumberFormat = "0.00%"
$ws.Range("C41").Value = 0.0267

# Row 42
$ws.Range("A42").Value = 0.11
$ws.Range("B42").NumberFormat = "0%"
$ws.Range("B42").Value = 0.95
$ws.Range("D42").Value = 18870
$ws.Range("C42").NumberFormat = "0.00%"
$ws.Range("C42").Formula = "=(D42-D31)/D31"

# Update page setup (A4, portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to D44, matching the new "next empty row" position
$ws.Range("D44").Select()
